$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.731.05'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '1.959.18'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''243.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('E7').Value = '  +5.38%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.372'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.21%  '
$ws.Range('D10').Value = '''0.0810'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '''22.15'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.04%  '
$ws.Range('D13').Value = '2.244.74'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').Value = '''0.820'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '1.961.08'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').Value = '36.676.13'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').Value = '''69.68'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').Value = '0.0₃0859'
$ws.Range('E20').Value = '  +1.28%  '
$ws.Range('D21').Value = '''5.10'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.01%  '
$ws.Range('D22').Value = '''228.13'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '''2.39'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('E25').Value = '  +3.57%  '
$ws.Range('D26').Value = '''9.29'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('D27').Value = '''0.137'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +16.08%  '
$ws.Range('D28').Value = '''160.96'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').Value = '''19.36'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('E30').Value = '  +2.26%  '
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('D33').Value = '''0.0618'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').Value = '''6.26'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +6.07%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  +20.79%  '
$ws.Range('E38').Value = '  +4.22%  '
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('D40').Value = '''0.100'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.28%  '
$ws.Range('D41').Value = '''2.89'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('E42').Value = '  +3.50%  '
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('D44').Value = '''16.06'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('E45').Value = '  +2.35%  '
$ws.Range('D46').Value = '1.345.30'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').Value = '''87.41'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('D48').Value = '''7.13'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('D50').Value = '2.136.91'
$ws.Range('E50').Value = '  +1.98%  '
$ws.Range('D51').Value = '''43.50'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.64%  '
